$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Reproduce the exact shared-string insertion order observed in the target
# workbook by writing cells in this sequence:
#   1) C9  -> "onyle context with userId (owner)"   (typo variant)
#   2) C8  -> "only context with userId (owner)"    (corrected variant)
#   3) A6  -> "fetchUser"
#   4) A7  -> "cancelUser"
#   5) B6  -> "targetUserId, txt"
#   6) B7  -> "targetUserId, txt"  (re-uses string from step 5)
#   7) C6  -> "userId from current queue; e.g. txt = 5min; invites and dequeues user"
#   8) C7  -> "userId from current queue; e.g. txt = "sry lul"; deuqueues user and cancel meeting"
#   9) A10 -> "requestUpdateOwnerScreen"

$ws.Range("A8").Value = "openLobby"
$ws.Range("A9").Value = "closeLobby"

$ws.Range("C9").Value = "onyle context with userId (owner)"
$ws.Range("C10").Value = "onyle context with userId (owner)"

$ws.Range("C8").Value = "only context with userId (owner)"

$ws.Range("A6").Value = "fetchUser"
$ws.Range("A7").Value = "cancelUser"

$ws.Range("B6").Value = "targetUserId, txt"
$ws.Range("B7").Value = "targetUserId, txt"

$ws.Range("C6").Value = "userId from current queue; e.g. txt = 5min; invites and dequeues user"
$ws.Range("C7").Value = 'userId from current queue; e.g. txt = "sry lul"; deuqueues user and cancel meeting'

$ws.Range("A10").Value = "requestUpdateOwnerScreen"

# Column width changes
$ws.Columns.Item(1).ColumnWidth = 23.06640625
$ws.Columns.Item(3).ColumnWidth = 65.6640625
